$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New configuration rows appended to the mapping table in column B.
$ws.Range("B4").Value = "Customer Reference No"
$ws.Range("B15").Value = "Labels"
$ws.Range("B16").Value = "Project"

# Reflect the last-selected cell as saved in the authored workbook.
$ws.Activate()
$ws.Range("A18").Select() | Out-Null
